# Apply the "Modelling VGG16 + CNN oversampling" sprint-board update.
#
# Summary of the edit (per commit message: "menambahkan model VGG16 dengan
# menggunakan teknik preprocessing dataset berupa oversampling"):
#   - Row 9  (task 7.)  "Preprocessing Dataset  Part 2" becomes
#                       "Preprocessing Dataset  Part 2 (Balancing Class)"
#                       and gets dates 07-11-2021 / 07-11-2021 / 07-11-2021,
#                       status "Done".
#   - Row 10 (task 8.)  "Modelling CNN 3 dan evaluasi" becomes
#                       "Modelling CNN dengan Oversampling Data"
#                       with dates 08-11-2021 / 09-11-2021 / 09-11-2021,
#                       status "Done".
#   - Row 11 (task 9.)  "Modelling CNN 4 dan evaluasi" becomes
#                       "Modelling VGG16 dengan Oversampling Data"
#                       with dates 08-11-2021 / 09-11-2021 / 09-11-2021,
#                       status "Done".
#   - Row 12 (task 10.) E12 picks up the same "text" number format as the
#                       other date columns (no value change, it stays blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dates in this sheet are stored as plain text (format code "@"), not real
# Excel date serials, so force the text number format before assigning the
# date-like strings to avoid them being auto-converted to date serials.
$ws.Range("C9:E11").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"

# --- Task labels first (matches the order new values were authored in) ---
$ws.Range("B9").Value = "Preprocessing Dataset  Part 2 (Balancing Class)"
$ws.Range("B11").Value = "Modelling VGG16 dengan Oversampling Data"
$ws.Range("B10").Value = "Modelling CNN dengan Oversampling Data"

# --- Row 10 & 11 dates: 08-11-2021 (start) / 09-11-2021 (finish) ---------
$ws.Range("C10").Value = "08-11-2021"
$ws.Range("C11").Value = "08-11-2021"
$ws.Range("D10").Value = "09-11-2021"
$ws.Range("E10").Value = "09-11-2021"
$ws.Range("D11").Value = "09-11-2021"
$ws.Range("E11").Value = "09-11-2021"

# --- Row 9 dates: 07-11-2021 -------------------------------------------
$ws.Range("C9").Value = "07-11-2021"
$ws.Range("D9").Value = "07-11-2021"
$ws.Range("E9").Value = "07-11-2021"

# --- Status column: mark all three tasks as Done (green fill) -----------
$ws.Range("G9").Value = "Done"
$ws.Range("G9").Interior.Color = 5287936
$ws.Range("G10").Value = "Done"
$ws.Range("G10").Interior.Color = 5287936
$ws.Range("G11").Value = "Done"
$ws.Range("G11").Interior.Color = 5287936

# The author's last selection before saving landed on D18.
$ws.Range("D18").Select()
